# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '60.727.97'
$ws.Cells.Item(2, 5).Value = '  +0.01%  '
$ws.Cells.Item(3, 4).Value = '2.338.47'
$ws.Cells.Item(3, 5).Value = '  -0.67%  '
$ws.Cells.Item(4, 5).Value = '  -0.16%  '
$ws.Cells.Item(5, 4).Value = '''547.30'
$ws.Cells.Item(5, 5).Value = '  +0.02%  '
$ws.Cells.Item(6, 4).Value = '''131.46'
$ws.Cells.Item(6, 5).Value = '  -0.71%  '
$ws.Cells.Item(7, 5).Value = '  -0.15%  '
$ws.Cells.Item(8, 4).Value = '''0.578'
$ws.Cells.Item(8, 5).Value = '  -1.78%  '
$ws.Cells.Item(9, 4).Value = '2.338.51'
$ws.Cells.Item(9, 5).Value = '  -0.49%  '
$ws.Cells.Item(10, 4).Value = '''0.102'
$ws.Cells.Item(10, 5).Value = '  +0.66%  '
$ws.Cells.Item(11, 4).Value = '''5.50'
$ws.Cells.Item(11, 5).Value = '  +0.03%  '
$ws.Cells.Item(12, 4).Value = '''0.149'
$ws.Cells.Item(12, 5).Value = '  -0.46%  '
$ws.Cells.Item(13, 4).Value = '''0.337'
$ws.Cells.Item(13, 5).Value = '  +0.64%  '
$ws.Cells.Item(14, 4).Value = '''23.59'
$ws.Cells.Item(14, 5).Value = '  -1.67%  '
$ws.Cells.Item(15, 4).Value = '2.756.53'
$ws.Cells.Item(15, 5).Value = '  -0.74%  '
$ws.Cells.Item(16, 4).Value = '60.643.16'
$ws.Cells.Item(16, 5).Value = '  -0.20%  '
$ws.Cells.Item(18, 4).Value = '2.343.76'
$ws.Cells.Item(18, 5).Value = '  -0.55%  '
$ws.Cells.Item(19, 4).Value = '''10.65'
$ws.Cells.Item(19, 5).Value = '  -0.65%  '
$ws.Cells.Item(20, 4).Value = '''4.09'
$ws.Cells.Item(20, 5).Value = '  -2.22%  '
$ws.Cells.Item(21, 4).Value = '''315.11'
$ws.Cells.Item(21, 5).Value = '  -0.11%  '
$ws.Cells.Item(22, 4).Value = '''6.58'
$ws.Cells.Item(22, 5).Value = '  -4.30%  '
$ws.Cells.Item(23, 4).Value = '''0.999'
$ws.Cells.Item(23, 5).Value = '  +0.17%  '
$ws.Cells.Item(24, 4).Value = '''63.91'
$ws.Cells.Item(24, 5).Value = '  +0.77%  '
$ws.Cells.Item(25, 5).Value = '  +0.53%  '
$ws.Cells.Item(26, 4).Value = '''1.00'
$ws.Cells.Item(26, 5).Value = '  -0.06%  '
$ws.Cells.Item(27, 4).Value = '''7.92'
$ws.Cells.Item(27, 5).Value = '  -0.61%  '
$ws.Cells.Item(28, 4).Value = '''1.40'
$ws.Cells.Item(28, 5).Value = '  +2.97%  '
$ws.Cells.Item(29, 5).Value = '  +9.08%  '
$ws.Cells.Item(30, 4).Value = '''172.83'
$ws.Cells.Item(30, 5).Value = '  +0.33%  '
$ws.Cells.Item(31, 5).Value = '  -2.07%  '
$ws.Cells.Item(32, 4).Value = '0.0₃0734'
$ws.Cells.Item(32, 5).Value = '  +0.31%  '
$ws.Cells.Item(33, 4).Value = '''5.96'
$ws.Cells.Item(33, 5).Value = '  +0.95%  '
$ws.Cells.Item(34, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(34, 4).Value = '''0.385'
$ws.Cells.Item(34, 5).Value = '  +0.53%  '
$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(35, 4).Value = '''1.37'
$ws.Cells.Item(35, 5).Value = '  -4.16%  '
$ws.Cells.Item(36, 4).Value = '''17.99'
$ws.Cells.Item(36, 5).Value = '  -0.49%  '
$ws.Cells.Item(37, 5).Value = '  +0.00%  '
$ws.Cells.Item(38, 5).Value = '  +0.01%  '
$ws.Cells.Item(39, 4).Value = '''4.14'
$ws.Cells.Item(39, 5).Value = '  -0.57%  '
$ws.Cells.Item(40, 4).Value = '''327.91'
$ws.Cells.Item(40, 5).Value = '  +4.02%  '
$ws.Cells.Item(41, 2).Value = 'Stacks'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(41, 4).Value = '''1.54'
$ws.Cells.Item(41, 5).Value = '  +0.21%  '
$ws.Cells.Item(42, 2).Value = 'OKB'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(42, 4).Value = '''38.16'
$ws.Cells.Item(42, 5).Value = '  -0.20%  '
$ws.Cells.Item(43, 4).Value = '''137.07'
$ws.Cells.Item(43, 5).Value = '  -3.90%  '
$ws.Cells.Item(44, 5).Value = '  +0.81%  '
$ws.Cells.Item(45, 4).Value = '''0.0943'
$ws.Cells.Item(45, 5).Value = '  -1.35%  '
$ws.Cells.Item(46, 4).Value = '''19.22'
$ws.Cells.Item(46, 5).Value = '  -0.75%  '
$ws.Cells.Item(47, 4).Value = '''0.570'
$ws.Cells.Item(47, 5).Value = '  +1.34%  '
$ws.Cells.Item(48, 4).Value = '''0.0497'
$ws.Cells.Item(48, 5).Value = '  -0.42%  '
$ws.Cells.Item(49, 5).Value = '  +1.11%  '
$ws.Cells.Item(50, 4).Value = '0.0₆0220'
$ws.Cells.Item(50, 5).Value = '  +4.60%  '
$ws.Cells.Item(51, 4).Value = '''10.99'
$ws.Cells.Item(51, 5).Value = '  -0.53%  '
